$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the data from columns D:F left into C:E (removing the
# "id_parametro" column's data without touching the <cols> width
# definitions for columns A-C).
for ($r = 1; $r -le 4; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $fVal
}

# Clear out the now-empty column F
$ws.Range("F1:F4").ClearContents()

# Update the selection to match the target state
$ws.Range("D7").Select()
